$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.257.49"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "1.592.96"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'213.24"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").Value = "'0.503"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Value = "'18.95"
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "1.818.21"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.595.55"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'4.00"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").Value = "'63.94"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "26.258.92"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").Value = "'214.61"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'4.29"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").Value = "'9.05"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("E24").Value = "  -3.97%  "
$ws.Range("D25").Value = "'145.03"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D29").Value = "'15.13"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "'3.19"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").Value = "1.418.31"
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Value = "'0.580"
$ws.Range("E37").Value = "  -3.82%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").Value = "'0.822"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").Value = "'5.80"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'0.982"
$ws.Range("E42").Value = "  -5.72%  "
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "1.729.90"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'60.99"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "'87.05"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "'1.49"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").Value = "'0.0509"
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("D50").Value = "'0.0952"
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("E51").Value = "  -0.05%  "
